$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the template values that the new rows (4-16) should repeat,
# except for column A (new address) and column B (incrementing counter).
$templateRow = 3
$address = "Москва, Ленинские горы, Д. 1"

for ($r = 4; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $address
    $ws.Cells.Item($r, 2).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($templateRow, 3).Value()
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($templateRow, 4).Value()
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($templateRow, 5).Value()
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($templateRow, 6).Value()
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($templateRow, 7).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($templateRow, 8).Value()
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($templateRow, 9).Value()
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($templateRow, 10).Value()
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($templateRow, 11).Value()
}

# Column A now needs to fit the much longer Moscow address, so it no
# longer keeps its old "best fit" narrow width - widen it to a custom width.
$ws.Columns.Item(1).ColumnWidth = 30.21875

# Match the selection left behind in the saved file.
$ws.Range("G15").Select() | Out-Null
